$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A4 to the new text value
$ws.Range("A4").Value = "路由为空时回home.html"

# Set B4 to the date value (serial 42790 => 2017-02-24), keep existing date format style
$ws.Range("B4").Value = (Get-Date -Year 2017 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0).Date

# Update the selected cell/range in the sheet view
$ws.Range("D13").Select()
